# Commit: "add word in june 7th"
#
# The last paragraph (empty) gets a new run "今天天气真好", formatted the
# same as the run in the paragraph right above it ("雨，最美不是下雨天，
# 而是和你一起躲过的屋檐"). The "_GoBack" bookmark, which currently sits
# at the end of that paragraph above, moves down to the end of the new
# run (Word keeps _GoBack pinned to the most recent edit position).

$d = $word.ActiveDocument

$paraCount = $d.Paragraphs.Count
$pSource = $d.Paragraphs.Item($paraCount - 1)   # "雨，最美……屋檐" paragraph
$pTarget = $d.Paragraphs.Item($paraCount)       # trailing empty paragraph

# 1) Re-home the "_GoBack" bookmark onto the (currently empty) last
#    paragraph first, collapsed at its end. Adding a bookmark with the
#    same name automatically removes/moves the previous one, so this
#    takes care of deleting it from the "雨，最美" paragraph too.
$bmAnchor = $pTarget.Range.Duplicate
$bmAnchor.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmAnchor)

# 2) Insert a placeholder run immediately in front of the bookmark, via
#    the bookmark's own Range -- InsertBefore keeps the new run ordered
#    before bookmarkStart/bookmarkEnd in the XML (run, then the two
#    bookmark markers), matching the target layout.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Range.InsertBefore("placeholder")

# 3) Swap the placeholder's formatting+text for an exact copy of the
#    source run's formatted text (rFonts hint="eastAsia", lang, etc.)
#    so the new run's rPr matches the sibling paragraph exactly.
$srcRange = $pSource.Range.Duplicate
$srcRange.End = $srcRange.End - 1   # exclude the paragraph mark

$placeholderRange = $d.Content
[void]$placeholderRange.Find.Execute("placeholder")
$placeholderRange.FormattedText = $srcRange.FormattedText

# 4) Fix up the text itself (the copy above still reads "雨，最美……屋檐").
$pTarget = $d.Paragraphs.Item($paraCount)
$textRange = $pTarget.Range.Duplicate
$textRange.End = $textRange.End - 1   # exclude the paragraph mark
$textRange.Text = "今天天气真好"
